$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3363.125
$ws.Range("I113").Value = 2801.6667
$ws.Range("J113").Value = 3700
$ws.Range("K113").Value = 2801.6667
$ws.Range("L113").Value = 3700
$ws.Range("M113").Value = 452.3332999999998
$ws.Range("N113").Value = -10208

$ws.Range("H125").Value = 1250
$ws.Range("J125").Value = 2000
$ws.Range("L125").Value = 18000
$ws.Range("N125").Value = -22920

$ws.Range("H135").Value = 622.7586
$ws.Range("J135").Value = 1675.875
$ws.Range("L135").Value = 15082.875
$ws.Range("N135").Value = -20152.875

$ws.Range("H137").Value = 1586.5
$ws.Range("I137").Value = 1604
$ws.Range("J137").Value = 1574
$ws.Range("K137").Value = 4812
$ws.Range("L137").Value = 4722
$ws.Range("M137").Value = -2262
$ws.Range("N137").Value = -9822

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4056.8408
$ws.Range("I32").Value = 3864.1892
$ws.Range("K32").Value = 3864.1892
$ws.Range("M32").Value = -3577.1892

$ws.Range("H61").Value = 1800
$ws.Range("I61").Value = 1100
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 1100
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -888
$ws.Range("N61").Value = -2924

$ws.Range("H74").Value = 1039.7
$ws.Range("I74").Value = 771
$ws.Range("J74").Value = 1666.6666
$ws.Range("K74").Value = 771
$ws.Range("L74").Value = 1666.6666
$ws.Range("M74").Value = 103
$ws.Range("N74").Value = -3414.6666

$ws.Range("H77").Value = 1039.7
$ws.Range("I77").Value = 771
$ws.Range("J77").Value = 1666.6666
$ws.Range("K77").Value = 3855
$ws.Range("L77").Value = 8333.333
$ws.Range("M77").Value = 513
$ws.Range("N77").Value = -17069.333

$ws.Range("H132").Value = 2617.762
$ws.Range("I132").Value = 2304
$ws.Range("J132").Value = 3036.111
$ws.Range("K132").Value = 6912
$ws.Range("L132").Value = 9108.332999999999
$ws.Range("M132").Value = -4382
$ws.Range("N132").Value = -14168.333

$ws.Range("H135").Value = 20024.166
$ws.Range("J135").Value = 20024.166
$ws.Range("L135").Value = 20024.166
$ws.Range("N135").Value = -30164.166

$ws.Range("H136").Value = 1800
$ws.Range("I136").Value = 1100
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 3300
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -750
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 249
$ws.Range("I7").Value = 248
$ws.Range("K7").Value = 248
$ws.Range("M7").Value = -135

$ws.Range("H134").Value = 9649.5625
$ws.Range("I134").Value = 925.3333
$ws.Range("K134").Value = 2775.9999
$ws.Range("M134").Value = -240.9998999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1796.5385
$ws.Range("I31").Value = 1317.2222
$ws.Range("K31").Value = 1317.2222
$ws.Range("M31").Value = -1022.2222

$ws.Range("H34").Value = 1796.5385
$ws.Range("I34").Value = 1317.2222
$ws.Range("K34").Value = 1317.2222
$ws.Range("M34").Value = -1115.2222

$ws.Range("H58").Value = 1744.6666
$ws.Range("I58").Value = 1515.5555
$ws.Range("J58").Value = 1973.7778
$ws.Range("K58").Value = 1515.5555
$ws.Range("L58").Value = 1973.7778
$ws.Range("M58").Value = -1312.5555
$ws.Range("N58").Value = -2379.7778

$ws.Range("H99").Value = 1620
$ws.Range("I99").Value = 1572.8572
$ws.Range("K99").Value = 1572.8572
$ws.Range("M99").Value = -74.85719999999992

$ws.Range("H126").Value = 1620
$ws.Range("I126").Value = 1572.8572
$ws.Range("K126").Value = 4718.571599999999
$ws.Range("M126").Value = -2248.571599999999

$ws.Range("H132").Value = 7465.3887
$ws.Range("I132").Value = 8783
$ws.Range("J132").Value = 4039.6
$ws.Range("K132").Value = 26349
$ws.Range("L132").Value = 12118.8
$ws.Range("M132").Value = -23819
$ws.Range("N132").Value = -17178.8

$ws.Range("H134").Value = 2137.0476
$ws.Range("I134").Value = 2255.7334
$ws.Range("K134").Value = 6767.2002
$ws.Range("M134").Value = -4232.2002

$ws.Range("H136").Value = 1744.6666
$ws.Range("I136").Value = 1515.5555
$ws.Range("J136").Value = 1973.7778
$ws.Range("K136").Value = 4546.666499999999
$ws.Range("L136").Value = 5921.3334
$ws.Range("M136").Value = -1996.666499999999
$ws.Range("N136").Value = -11021.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 3200
$ws.Range("J87").Value = 3200
$ws.Range("L87").Value = 9600
$ws.Range("N87").Value = -12096

$ws.Range("H90").Value = 3200
$ws.Range("J90").Value = 3200
$ws.Range("L90").Value = 28800
$ws.Range("N90").Value = -41280

$ws.Range("H92").Value = 765.8889
$ws.Range("J92").Value = 800
$ws.Range("L92").Value = 2400
$ws.Range("N92").Value = -4896

$ws.Range("H132").Value = 1550
$ws.Range("I132").Value = 1200
$ws.Range("J132").Value = 1900
$ws.Range("K132").Value = 10800
$ws.Range("L132").Value = 17100
$ws.Range("M132").Value = -8270
$ws.Range("N132").Value = -22160

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2817.7896
$ws.Range("I132").Value = 2467.2144
$ws.Range("K132").Value = 7401.6432
$ws.Range("M132").Value = -4871.6432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2557.7144
$ws.Range("I7").Value = 2408.5833
$ws.Range("J7").Value = 3452.5
$ws.Range("K7").Value = 2408.5833
$ws.Range("L7").Value = 3452.5
$ws.Range("M7").Value = -2296.5833
$ws.Range("N7").Value = -3676.5

$ws.Range("H68").Value = 1722.8572
$ws.Range("I68").Value = 1635
$ws.Range("J68").Value = 2250
$ws.Range("K68").Value = 1635
$ws.Range("L68").Value = 2250
$ws.Range("M68").Value = -886
$ws.Range("N68").Value = -3748

$ws.Range("H71").Value = 1722.8572
$ws.Range("I71").Value = 1635
$ws.Range("J71").Value = 2250
$ws.Range("K71").Value = 8175
$ws.Range("L71").Value = 11250
$ws.Range("M71").Value = -4431
$ws.Range("N71").Value = -18738

$ws.Range("H93").Value = 2000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = ""
$ws.Range("N93").Value = -4496

$ws.Range("H126").Value = 2557.7144
$ws.Range("I126").Value = 2408.5833
$ws.Range("J126").Value = 3452.5
$ws.Range("K126").Value = 7225.749899999999
$ws.Range("L126").Value = 10357.5
$ws.Range("M126").Value = -4755.749899999999
$ws.Range("N126").Value = -15297.5

$ws.Range("H132").Value = 20830.596
$ws.Range("I132").Value = 1323.2413
$ws.Range("J132").Value = 45426.824
$ws.Range("K132").Value = 3969.7239
$ws.Range("L132").Value = 136280.472
$ws.Range("M132").Value = -1439.7239
$ws.Range("N132").Value = -141340.472

$ws.Range("H136").Value = 4252.25
$ws.Range("I136").Value = 5122.72
$ws.Range("K136").Value = 15368.16
$ws.Range("M136").Value = -12818.16

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2326.3823
$ws.Range("I132").Value = 1819.4231
$ws.Range("J132").Value = 3974
$ws.Range("K132").Value = 5458.2693
$ws.Range("L132").Value = 11922
$ws.Range("M132").Value = -2928.2693
$ws.Range("N132").Value = -16982

$ws.Range("H136").Value = 760.4
$ws.Range("I136").Value = 755
$ws.Range("J136").Value = 761.75
$ws.Range("K136").Value = 2265
$ws.Range("L136").Value = 2285.25
$ws.Range("M136").Value = 285
$ws.Range("N136").Value = -7385.25
